# Auto-generated Excel COM-interop script
# Applies the BRVM daily-refresh diff to 'Recommandations' and 'Top_YTD' sheets.

$wb = $excel.ActiveWorkbook
$wsReco = $wb.Worksheets.Item("Recommandations")
$wsYtd  = $wb.Worksheets.Item("Top_YTD")

# ---- Recommandations sheet ----
$wsReco.Cells.Item(2, 4).Value = 3730
$wsReco.Cells.Item(2, 5).Value = 900
$wsReco.Cells.Item(3, 4).Value = 3382.45
$wsReco.Cells.Item(3, 5).Value = 112.06
$wsReco.Cells.Item(4, 5).Value = 700
$wsReco.Cells.Item(5, 4).Value = 2364.16
$wsReco.Cells.Item(5, 5).Value = 590.14
$wsReco.Cells.Item(6, 4).Value = 2022.5
$wsReco.Cells.Item(6, 5).Value = 513.02
$wsReco.Cells.Item(7, 4).Value = 1441.71
$wsReco.Cells.Item(7, 5).Value = 354.92
$wsReco.Cells.Item(8, 4).Value = 1349.63
$wsReco.Cells.Item(8, 5).Value = 339.25
$wsReco.Cells.Item(10, 4).Value = 691.03
$wsReco.Cells.Item(10, 5).Value = 176.86
$wsReco.Cells.Item(11, 1).Value = 'BRVM - FINANCES'
$wsReco.Cells.Item(11, 3).Value = 4
$wsReco.Cells.Item(11, 4).Value = 582.08
$wsReco.Cells.Item(11, 5).Value = 144.93
$wsReco.Cells.Item(12, 1).Value = 'BRVM - SERVICES FINANCIERS'
$wsReco.Cells.Item(12, 4).Value = 572.06
$wsReco.Cells.Item(12, 5).Value = 142.44
$wsReco.Cells.Item(13, 1).Value = 'BRVM-PRESTIGE'
$wsReco.Cells.Item(13, 4).Value = 564.33
$wsReco.Cells.Item(13, 5).Value = 141.39
$wsReco.Cells.Item(14, 1).Value = 'BRVM - INDUSTRIE                 (**)'
$wsReco.Cells.Item(14, 3).Value = 2
$wsReco.Cells.Item(14, 4).Value = 527.65
$wsReco.Cells.Item(14, 5).Value = 266.85
$wsReco.Cells.Item(15, 1).Value = 'BRVM - INDUSTRIELS'
$wsReco.Cells.Item(15, 3).Value = 4
$wsReco.Cells.Item(15, 4).Value = 491.13
$wsReco.Cells.Item(15, 5).Value = 121.16
$wsReco.Cells.Item(16, 1).Value = 'BRVM - ENERGIE'
$wsReco.Cells.Item(16, 4).Value = 451.62
$wsReco.Cells.Item(16, 5).Value = 111.8
$wsReco.Cells.Item(17, 1).Value = 'BRVM-PRINCIPAL                    (**)'
$wsReco.Cells.Item(17, 3).Value = 2
$wsReco.Cells.Item(17, 4).Value = 438.39
$wsReco.Cells.Item(17, 5).Value = 219.41
$wsReco.Cells.Item(18, 1).Value = 'BRVM - TELECOMMUNICATIONS'
$wsReco.Cells.Item(18, 3).Value = 4
$wsReco.Cells.Item(18, 4).Value = 380.61
$wsReco.Cells.Item(18, 5).Value = 96.88
$wsReco.Cells.Item(19, 1).Value = 'BRVM - CONSOMMATION DE BASE          (**)'
$wsReco.Cells.Item(19, 3).Value = 1
$wsReco.Cells.Item(19, 4).Value = 222.65
$wsReco.Cells.Item(19, 5).Value = 222.65
$wsReco.Cells.Item(20, 1).Value = 'BRVM - CONSOMMATION DE BASE         (**)'
$wsReco.Cells.Item(20, 4).Value = 217.81
$wsReco.Cells.Item(20, 5).Value = 217.81
$wsReco.Cells.Item(21, 1).Value = 'UNILEVER CI (UNLC)'
$wsReco.Cells.Item(21, 2).Value = 2
$wsReco.Cells.Item(21, 3).Value = 0
$wsReco.Cells.Item(21, 4).Value = 14.97
$wsReco.Cells.Item(21, 5).Value = 7.47
$wsReco.Cells.Item(22, 1).Value = 'UNIWAX CI (UNXC)'
$wsReco.Cells.Item(22, 4).Value = 9.94
$wsReco.Cells.Item(22, 5).Value = 7.14
$wsReco.Cells.Item(23, 1).Value = 'BANK OF AFRICA ML (BOAM)'
$wsReco.Cells.Item(23, 2).Value = 1
$wsReco.Cells.Item(23, 4).Value = 5.28
$wsReco.Cells.Item(23, 5).Value = 5.28
$wsReco.Cells.Item(23, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(23, 7).Value = '➖ Neutre'
$wsReco.Cells.Item(24, 4).Value = 4.97
$wsReco.Cells.Item(24, 5).Value = 1.32
$wsReco.Cells.Item(25, 1).Value = 'SAPH CI (SPHC)'
$wsReco.Cells.Item(25, 2).Value = 1
$wsReco.Cells.Item(25, 3).Value = 0
$wsReco.Cells.Item(25, 4).Value = 4.7
$wsReco.Cells.Item(25, 5).Value = 4.7
$wsReco.Cells.Item(25, 7).Value = '➖ Neutre'
$wsReco.Cells.Item(26, 1).Value = 'ORANGE COTE D''IVOIRE (ORAC)'
$wsReco.Cells.Item(26, 4).Value = 4.3
$wsReco.Cells.Item(26, 5).Value = 4.3
$wsReco.Cells.Item(27, 1).Value = 'BICI CI (BICC)'
$wsReco.Cells.Item(27, 4).Value = 3.85
$wsReco.Cells.Item(27, 5).Value = 3.85
$wsReco.Cells.Item(28, 1).Value = 'NESTLE CI (NTLC)'
$wsReco.Cells.Item(28, 4).Value = 3.04
$wsReco.Cells.Item(28, 5).Value = 3.04
$wsReco.Cells.Item(29, 1).Value = 'CFAO MOTORS CI (CFAC)'
$wsReco.Cells.Item(29, 4).Value = 2.67
$wsReco.Cells.Item(29, 5).Value = 4.98
$wsReco.Cells.Item(32, 1).Value = 'BANK OF AFRICA SENEGAL (BOAS)'
$wsReco.Cells.Item(32, 2).Value = 1
$wsReco.Cells.Item(32, 3).Value = 1
$wsReco.Cells.Item(32, 4).Value = 0.13
$wsReco.Cells.Item(32, 5).Value = 3.58
$wsReco.Cells.Item(32, 7).Value = '👀 À surveiller'
$wsReco.Cells.Item(33, 1).Value = 'TOTAL'
$wsReco.Cells.Item(33, 2).Value = 0
$wsReco.Cells.Item(33, 4).Value = 0
$wsReco.Cells.Item(33, 5).Value = 0
$wsReco.Cells.Item(33, 7).Value = '➖ Neutre'
$wsReco.Cells.Item(35, 1).Value = 'SAFCA CI (SAFC)'
$wsReco.Cells.Item(35, 3).Value = 1
$wsReco.Cells.Item(35, 4).Value = -0.56
$wsReco.Cells.Item(35, 5).Value = 5.02
$wsReco.Cells.Item(36, 1).Value = 'SOLIBRA CI (SLBC)'
$wsReco.Cells.Item(36, 2).Value = 1
$wsReco.Cells.Item(36, 4).Value = -1.58
$wsReco.Cells.Item(36, 5).Value = 3.67
$wsReco.Cells.Item(36, 7).Value = '👀 À surveiller'
$wsReco.Cells.Item(37, 1).Value = 'AFRICA GLOBAL LOGISTICS CI (SDSC)'
$wsReco.Cells.Item(37, 4).Value = -1.67
$wsReco.Cells.Item(37, 5).Value = -1.67
$wsReco.Cells.Item(38, 1).Value = 'SERVAIR ABIDJAN CI (ABJC)'
$wsReco.Cells.Item(38, 4).Value = -2.5
$wsReco.Cells.Item(38, 5).Value = -2.5
$wsReco.Cells.Item(39, 1).Value = 'SMB CI (SMBC)'
$wsReco.Cells.Item(39, 4).Value = -3.11
$wsReco.Cells.Item(39, 5).Value = -3.11
$wsReco.Cells.Item(42, 1).Value = 'NEI-CEDA CI (NEIC)'
$wsReco.Cells.Item(42, 2).Value = 0
$wsReco.Cells.Item(42, 4).Value = -4.26
$wsReco.Cells.Item(42, 5).Value = -4.26
$wsReco.Cells.Item(42, 7).Value = '➖ Neutre'
$wsReco.Cells.Item(45, 1).Value = 'SICOR CI (SICC)'
$wsReco.Cells.Item(45, 4).Value = -6.87
$wsReco.Cells.Item(45, 5).Value = -6.87
$wsReco.Cells.Item(46, 1).Value = 'SICABLE CI (CABC)'
$wsReco.Cells.Item(46, 2).Value = 1
$wsReco.Cells.Item(46, 3).Value = 2
$wsReco.Cells.Item(46, 4).Value = -7.34
$wsReco.Cells.Item(46, 5).Value = -7.4
$wsReco.Cells.Item(46, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(46, 7).Value = '👀 À surveiller'

# ---- Top_YTD sheet ----
$wsYtd.Cells.Item(2, 2).Value = 10027287.37
$wsYtd.Cells.Item(3, 2).Value = 1135606
$wsYtd.Cells.Item(5, 2).Value = 227894.29
$wsYtd.Cells.Item(6, 2).Value = 134415
$wsYtd.Cells.Item(7, 2).Value = 44835.98
$wsYtd.Cells.Item(8, 2).Value = 36502.7
$wsYtd.Cells.Item(9, 2).Value = 5433.87
$wsYtd.Cells.Item(10, 2).Value = 3533.67
$wsYtd.Cells.Item(11, 2).Value = 3387.63

